$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch2")
$ws.Rows.Item(6).Insert()
$ws.Range("B6").Value = "The procedures described in this and later chapters require the Internetworking environment, please connect your computer to the network. If you don't have the network, please refer the text for the environment without the Internetworking."
Write-Host "done"
